$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowFV {
    param($row, $vals)
    $col = 6
    foreach ($v in $vals) {
        $ws.Cells.Item($row, $col).Value = $v
        $col++
    }
}

# Rows 13 and 14 swap their match data (F:V)
Set-RowFV 13 @("Hillerod", 0, "Hobro", 1, 2.23, "29/07/2023 14:12", 3.15, "04/08/2023 18:23", 3.38, "29/07/2023 14:12", 3.49, "04/08/2023 18:23", 3.23, "29/07/2023 14:12", 2.27, "04/08/2023 18:23", "https://www.betexplorer.com/football/denmark/1st-division/hillerod-hobro/E3bkjiAO/")
Set-RowFV 14 @("Horsens", 2, "Koge", 0, 1.72, "01/08/2023 05:42", 1.81, "04/08/2023 18:52", 3.95, "01/08/2023 05:42", 4.16, "04/08/2023 18:56", 4.12, "01/08/2023 05:42", 3.96, "04/08/2023 18:56", "https://www.betexplorer.com/football/denmark/1st-division/horsens-koge/zDaoiXfI/")

# Rows 43, 44, 45 rotate: new43=old45, new44=old43, new45=old44
Set-RowFV 43 @("Vendsyssel", 0, "Fredericia", 2, 2.02, "27/08/2023 13:12", 2.46, "01/09/2023 18:55", 3.79, "27/08/2023 13:12", 3.9, "01/09/2023 18:55", 3.19, "27/08/2023 13:12", 2.62, "01/09/2023 18:55", "https://www.betexplorer.com/football/denmark/1st-division/vendsyssel-ff-fredericia/zmIxmLb8/")
Set-RowFV 44 @("Horsens", 0, "B.93", 0, 1.79, "28/08/2023 01:42", 1.59, "01/09/2023 18:57", 4.02, "28/08/2023 01:42", 4.54, "01/09/2023 18:57", 3.72, "28/08/2023 01:42", 4.97, "01/09/2023 18:57", "https://www.betexplorer.com/football/denmark/1st-division/horsens-boldklubben-1893/MqUfOyM7/")
Set-RowFV 45 @("Hillerod", 2, "Sonderjyske", 2, 3.9, "28/08/2023 18:42", 4.16, "01/09/2023 18:58", 3.8, "28/08/2023 18:42", 3.85, "01/09/2023 18:58", 1.79, "28/08/2023 18:42", 1.83, "01/09/2023 18:58", "https://www.betexplorer.com/football/denmark/1st-division/hillerod-sonderjyske/EwHtnuEE/")

# Rows 68 and 69 swap their match data (F:V)
Set-RowFV 68 @("Horsens", 1, "Aalborg", 4, 3.66, "29/09/2023 18:13", 3.88, "06/10/2023 18:53", 3.86, "29/09/2023 18:13", 3.72, "06/10/2023 18:54", 1.84, "29/09/2023 18:13", 1.93, "06/10/2023 18:54", "https://www.betexplorer.com/football/denmark/1st-division/horsens-aalborg/SUv8VC77/")
Set-RowFV 69 @("B.93", 0, "Hillerod", 3, 2.6, "30/09/2023 13:12", 3.6, "06/10/2023 18:59", 3.5, "30/09/2023 13:12", 3.81, "06/10/2023 18:59", 2.5, "30/09/2023 13:12", 1.98, "06/10/2023 18:59", "https://www.betexplorer.com/football/denmark/1st-division/boldklubben-1893-hillerod/hCneYENl/")

# New rows 79-84 appended; copy formatting for columns A and E from row 2
foreach ($r in 79..84) {
    $ws.Range("A2").Copy()
    $ws.Range("A$r").PasteSpecial(-4122)
    $ws.Range("E2").Copy()
    $ws.Range("E$r").PasteSpecial(-4122)
}

$ws.Cells.Item(79, 1).Value = 78
$ws.Cells.Item(79, 2).Value = "denmark"
$ws.Cells.Item(79, 3).Value = "1st-division"
$ws.Cells.Item(79, 4).Value = "2023-2024"
$ws.Cells.Item(79, 5).Value = 45226.79166666666
Set-RowFV 79 @("Hobro", 2, "Koge", 1, 1.74, "22/10/2023 16:12", 1.65, "27/10/2023 18:51", 3.92, "22/10/2023 16:12", 4.26, "27/10/2023 18:51", 4.42, "22/10/2023 16:12", 4.85, "27/10/2023 18:51", "https://www.betexplorer.com/football/denmark/1st-division/hobro-koge/UBGxoh7f/")

$ws.Cells.Item(80, 1).Value = 79
$ws.Cells.Item(80, 2).Value = "denmark"
$ws.Cells.Item(80, 3).Value = "1st-division"
$ws.Cells.Item(80, 4).Value = "2023-2024"
$ws.Cells.Item(80, 5).Value = 45226.79166666666
Set-RowFV 80 @("B.93", 0, "Sonderjyske", 4, 5.01, "22/10/2023 15:12", 8.03, "27/10/2023 18:58", 4.47, "22/10/2023 15:12", 5.46, "27/10/2023 18:58", 1.57, "22/10/2023 15:12", 1.34, "27/10/2023 18:50", "https://www.betexplorer.com/football/denmark/1st-division/boldklubben-1893-sonderjyske/2PEpqWy7/")

$ws.Cells.Item(81, 1).Value = 80
$ws.Cells.Item(81, 2).Value = "denmark"
$ws.Cells.Item(81, 3).Value = "1st-division"
$ws.Cells.Item(81, 4).Value = "2023-2024"
$ws.Cells.Item(81, 5).Value = 45226.79166666666
Set-RowFV 81 @("Horsens", 1, "Hillerod", 3, 1.84, "20/10/2023 19:12", 2.17, "27/10/2023 18:51", 3.9, "20/10/2023 19:12", 3.59, "27/10/2023 18:51", 3.63, "20/10/2023 19:12", 3.27, "27/10/2023 18:51", "https://www.betexplorer.com/football/denmark/1st-division/horsens-hillerod/G4dX6ZTD/")

$ws.Cells.Item(82, 1).Value = 81
$ws.Cells.Item(82, 2).Value = "denmark"
$ws.Cells.Item(82, 3).Value = "1st-division"
$ws.Cells.Item(82, 4).Value = "2023-2024"
$ws.Cells.Item(82, 5).Value = 45227.54166666666
Set-RowFV 82 @("Kolding IF", 0, "Aalborg", 1, 3.31, "22/10/2023 15:12", 3.64, "28/10/2023 12:57", 3.54, "22/10/2023 15:12", 3.66, "28/10/2023 12:58", 2.13, "22/10/2023 15:12", 2.01, "28/10/2023 12:57", "https://www.betexplorer.com/football/denmark/1st-division/kolding-if-aalborg/hUJ0k7bt/")

$ws.Cells.Item(83, 1).Value = 82
$ws.Cells.Item(83, 2).Value = "denmark"
$ws.Cells.Item(83, 3).Value = "1st-division"
$ws.Cells.Item(83, 4).Value = "2023-2024"
$ws.Cells.Item(83, 5).Value = 45227.75
Set-RowFV 83 @("Helsingor", 3, "Fredericia", 3, 2.97, "21/10/2023 18:13", 3.93, "28/10/2023 17:42", 3.69, "21/10/2023 18:13", 3.94, "28/10/2023 17:42", 2.15, "21/10/2023 18:13", 1.85, "28/10/2023 17:42", "https://www.betexplorer.com/football/denmark/1st-division/helsingor-fredericia/EaIlrjiD/")

$ws.Cells.Item(84, 1).Value = 83
$ws.Cells.Item(84, 2).Value = "denmark"
$ws.Cells.Item(84, 3).Value = "1st-division"
$ws.Cells.Item(84, 4).Value = "2023-2024"
$ws.Cells.Item(84, 5).Value = 45228.54166666666
Set-RowFV 84 @("Naestved", 1, "Vendsyssel", 1, 2.46, "23/10/2023 09:12", 2.61, "29/10/2023 12:43", 3.47, "23/10/2023 09:12", 3.91, "29/10/2023 12:43", 2.66, "23/10/2023 09:12", 2.47, "29/10/2023 12:43", "https://www.betexplorer.com/football/denmark/1st-division/naestved-if-vendsyssel-ff/82FtpCM0/")

